$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 717 (shifts old rows 717-758 down to 719-760),
# matching cell formatting/format of the row being pushed down.
$ws.Rows.Item(717).Insert()
$ws.Rows.Item(717).Insert()

# Populate the two newly inserted rows with the new data points.
# The date column stores plain text like "2026/01/24" (not a real Excel date),
# so prefix with an apostrophe to force text and avoid Excel's date auto-parse.
$ws.Range("A717").Value = "'2026/01/24"
$ws.Range("A717").Style = "Normal"
$ws.Range("B717").Value = "土"
$ws.Range("C717").Value = 22
$ws.Range("D717").Value = 158

$ws.Range("A718").Value = "'2026/01/25"
$ws.Range("A718").Style = "Normal"
$ws.Range("B718").Value = "日"
$ws.Range("C718").Value = 2
$ws.Range("D718").Value = 167
